# Add a new "JAR" entry (command reference + follow-up note) and a new
# "i18n" entry to the Java reference sheet, appended as rows 125-127
# below the existing last row (124). Mirrors the commit:
# "Add java entry about i18n(internationalization)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$jarTitle = @'
JAR
'@

$jarBody = @'
# JAR command
Usage: jar {ctxui}[vfmn0PMe] [jar-file] [manifest-file] [entry-point] [-C dir] files ...
Options:
    -c  create new archive
    -t  list table of contents for archive
    -x  extract named (or all) files from archive
    -u  update existing archive
    -v  generate verbose output on standard output
    -f  specify archive file name
    -m  include manifest information from specified manifest file
    -n  perform Pack200 normalization after creating a new archive
    -e  specify application entry point for stand-alone application
        bundled into an executable jar file
    -0  store only; use no ZIP compression
    -P  preserve leading '/' (absolute path) and ".." (parent directory) components from file names
    -M  do not create a manifest file for the entries
    -i  generate index information for the specified jar files
    -C  change to the specified directory and include the following file
If any file is a directory then it is processed recursively.
The manifest file name, the archive file name and the entry point name are
specified in the same order as the 'm', 'f' and 'e' flags.
Example 1: to archive two class files into an archive called classes.jar:
       jar cvf classes.jar Foo.class Bar.class
Example 2: use an existing manifest file 'mymanifest' and archive all the
           files in the foo/ directory into 'classes.jar':
       jar cvfm classes.jar mymanifest -C foo/ .
'@

$jar2Title = @'
JAR: more on jar
'@

$jar2Body = @'
We can use "jar" command to add info to manifest file from another file or from command line args. But with a build tool like ant or maven, this is replace by cnetralized project doc like "pom.xml"
'@

$i18nTitle = @'
i18n
'@

$i18nBody = @'
i18n is the abbreviation of internationalization. Similiarly, l10n is the abbreviation of localization. Util ResourceBundle and Locale are provided in JDK to read properties from bundle file in a certainly formating like "MessagesBundle_fr_FR.properties" and "MessagesBundle_en_US.properties" with code below:
```
import java.util.*;
public class I18NSample {
    static public void main(String[] args) {
        String language;
        String country;
        if (args.length != 2) {
              return; //invalid number of ags
        } else {
            language = new String(args[0]);
            country = new String(args[1]);
        }
        Locale currentLocale = new Locale(language, country);
        ResourceBundle messages = ResourceBundle.getBundle("MessagesBundle", currentLocale);
        System.out.println(messages.getString("greetings"));
        System.out.println(messages.getString("inquiry"));
        System.out.println(messages.getString("farewell"));
    }
}
```
While properties file follow classic format:
```
greetings = Bonjour.
farewell = Au revoir.
inquiry = Comment allez-vous?
```
'@

# Row 125: JAR command reference (keeps the "RTFM" category used by the
# surrounding rows). Insert a copy of the last existing row (124) so the
# new row inherits its style/height, then overwrite its values.
$ws.Rows.Item(124).Copy()
$ws.Rows.Item(125).Insert(-4121)
$ws.Range("A125").Value = "RTFM"
$ws.Range("B125").Value = $jarTitle
$ws.Range("C125").Value = $jarBody
$ws.Rows.Item(125).RowHeight = 32.25

# Row 126: follow-up note on "jar" / manifest files.
$ws.Rows.Item(125).Copy()
$ws.Rows.Item(126).Insert(-4121)
$ws.Range("A126").Value = "RTFM"
$ws.Range("B126").Value = $jar2Title
$ws.Range("C126").Value = $jar2Body
$ws.Rows.Item(126).RowHeight = 32.25

# Row 127: new i18n (internationalization) entry.
$ws.Rows.Item(126).Copy()
$ws.Rows.Item(127).Insert(-4121)
$ws.Range("A127").Value = "RTFM"
$ws.Range("B127").Value = $i18nTitle
$ws.Range("C127").Value = $i18nBody
$ws.Rows.Item(127).RowHeight = 32.25

$ws.Application.CutCopyMode = $false

# Move the view/selection the way the author left it after typing the
# new rows: scrolled down a bit further, cursor parked a couple of rows
# below the newly-added data.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 117
$win.ScrollColumn = 1
$ws.Range("A129").Select()
